$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: account holder name / card number change
$ws.Range("C2").Value = "Hartmut"
# Card number is a long digit string that must stay text (not be coerced to a
# number, which would lose its exact display as a 16-digit card number).
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 09.11.2024"

# Row 6
$ws.Range("B6").Value = "11.11."
$ws.Range("C6").Value = "12.11."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 93326309"
$ws.Range("E6").Value = "85,87-"

# Row 7
$ws.Range("B7").Value = "15.11."
$ws.Range("C7").Value = "16.11."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-29240970"
$ws.Range("E7").Value = "55,55-"

# Row 8
$ws.Range("B8").Value = "19.11."
$ws.Range("C8").Value = "20.11."
$ws.Range("D8").Value = "EBAY MKTPLC EU YOEFUK"
$ws.Range("E8").Value = "148,75-"

# Rows 9-11: clear remaining transactions (now blank)
$ws.Range("B9").Value = $null
$ws.Range("C9").Value = $null
$ws.Range("D9").Value = $null
$ws.Range("E9").Value = $null
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

$ws.Range("B10").Value = $null
$ws.Range("C10").Value = $null
$ws.Range("D10").Value = $null
$ws.Range("E10").Value = $null
$ws.Range("E10").HorizontalAlignment = -4152
$ws.Range("E10").VerticalAlignment = -4108
$ws.Range("E10").WrapText = $true

$ws.Range("B11").Value = $null
$ws.Range("C11").Value = $null
$ws.Range("D11").Value = $null
$ws.Range("E11").Value = $null
$ws.Range("E11").HorizontalAlignment = -4152
$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("E11").WrapText = $true

# Row 12: closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 22.11.2024"
$ws.Range("E12").Value = "290,17-"

# Row 13: next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 02.12.2024"
